$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing activity rows (2-7) down by two rows (to 4-9),
# working from the bottom up so we don't overwrite data we still need.
for ($r = 7; $r -ge 2; $r--) {
    $destRow = $r + 2
    for ($c = 1; $c -le 9; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $ws.Cells.Item($r, $c).Value2
    }
}

# New row 2: Lunch Run
$ws.Cells.Item(2, 1).Value2 = "Lunch Run"
$ws.Cells.Item(2, 2).Value2 = 12
$ws.Cells.Item(2, 3).Value2 = "01:18:08"
$ws.Cells.Item(2, 4).Value2 = 4688
$ws.Cells.Item(2, 5).Value2 = "Run"
$ws.Cells.Item(2, 6).Value2 = "2025-07-09T12:24:49Z"
$ws.Cells.Item(2, 7).Value2 = "06:31"
$ws.Cells.Item(2, 8).Value2 = 9.752624640000001
$ws.Cells.Item(2, 9).Value2 = 139.5

# New row 3: Säbä
$ws.Cells.Item(3, 1).Value2 = "Säbä"
$ws.Cells.Item(3, 2).Value2 = 6.41
$ws.Cells.Item(3, 3).Value2 = "01:38:20"
$ws.Cells.Item(3, 4).Value2 = 5900
$ws.Cells.Item(3, 5).Value2 = "Run"
$ws.Cells.Item(3, 6).Value2 = "2025-07-08T19:33:57Z"
$ws.Cells.Item(3, 7).Value2 = "15:20"
$ws.Cells.Item(3, 8).Value2 = 10.68604416
$ws.Cells.Item(3, 9).Value2 = 131.2
